# "Correct the weird bug": a new "Date of Last Update" column is inserted
# after "Time Worked (Minutes)", and an extra (duplicate) "TicketID" column
# is inserted right after it, pushing the two trailing columns
# ("Quality & Continuous Improvement / ACE" / "Quality-Related Issues") two
# places to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new columns starting at column I -------------------------
# 1st insert makes room for the new "Date of Last Update" column.
$ws.Columns("I").Insert()
# 2nd insert makes room for the duplicated "TicketID" column.
$ws.Columns("J").Insert()

# --- Fill in the two new columns ------------------------------------------
$ws.Range("I1").Value = "Date of Last Update"
$ws.Range("I2").Value = "2017-08-15T09:28:22.373000"

$ws.Range("J1").Value = "TicketID"
# Format as Text first so the numeric-looking ticket id "249" is stored as a
# string (matching column A), not auto-converted to a number, then drop back
# to the Normal style so no extra number format lingers on the cell.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "249"
$ws.Range("J2").Style = "Normal"

# --- Fix up the value that landed in the wrong row after the shift --------
$ws.Range("K2").Value = "Quality-Related Issues"

# --- Match formatting of the header row (bold) across the new columns -----
$ws.Range("A1:L1").Font.Bold = $true
$ws.Range("A2:L2").Font.Bold = $false

# --- Misc cosmetic bits from the diff --------------------------------------
$null = $ws.Range("A1").Select()
